$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column J/K previously held shared-string placeholder text ("r"/"s") on
# row 1 only; the data source was corrected so row 1 now carries real
# numeric values like every other row.
$ws.Range("J1").Value = 0.3
$ws.Range("K1").Value = 1

# The whole K column (rows 2-51) was updated from the old constant 0.3
# to the corrected constant 1.
$ws.Range("K2:K51").Value = 1

# Reflect the updated view/selection left by the author: window scrolled
# one row further down and the live selection now sits on K1:K51.
$win = $excel.ActiveWindow
$win.ScrollRow = 39
$win.ScrollColumn = 1
$ws.Range("K1:K51").Select() | Out-Null
